# Auto-generated Excel COM-interop script to apply the 'Fix heat rate modeling syntax' edit
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("DG Dispatch")
$ws.Range("S2").Value = 0.0

$ws = $wb.Worksheets.Item("Costs and Revenues")
$ws.Range("B2").Value = 76271.0624
$ws.Range("D2").Value = 9300.638068405267
$ws.Range("F2").Value = 50991.37327170695

$ws = $wb.Worksheets.Item("PV Dispatch")
$ws.Range("H2").Value = 41.6
$ws.Range("I2").Value = 50.43636363636368
$ws.Range("K3").Value = 0.0
$ws.Range("M3").Value = 104.0
$ws.Range("N3").Value = 47.45311702887471
$ws.Range("Q3").Value = 52.0
$ws.Range("J4").Value = 0.0
$ws.Range("K4").Value = 0.0
$ws.Range("L4").Value = 0.0
$ws.Range("M4").Value = 83.2
$ws.Range("N4").Value = 83.2
$ws.Range("P4").Value = 29.58312417100299
$ws.Range("Q4").Value = 0.0

$ws = $wb.Worksheets.Item("Battery Input")
$ws.Range("H2").Value = 28.6
$ws.Range("I2").Value = 19.23636363636368
$ws.Range("K3").Value = 0.0
$ws.Range("M3").Value = 80.6
$ws.Range("N3").Value = 21.45311702887471
$ws.Range("Q3").Value = 26.0
$ws.Range("J4").Value = 0.0
$ws.Range("K4").Value = 0.0
$ws.Range("L4").Value = 0.0
$ws.Range("M4").Value = 59.8
$ws.Range("N4").Value = 83.2
$ws.Range("P4").Value = 29.58312417100299
$ws.Range("Q4").Value = 0.0

$ws = $wb.Worksheets.Item("Battery Output")
$ws.Range("S2").Value = 7.515999999999954

$ws = $wb.Worksheets.Item("State of Charge")
$ws.Range("B2").Value = 188.6909090909091
$ws.Range("C2").Value = 168.9939393939394
$ws.Range("D2").Value = 155.8626262626263
$ws.Range("E2").Value = 142.7313131313131
$ws.Range("G2").Value = 142.47
$ws.Range("H2").Value = 170.784
$ws.Range("S2").Value = 640.4080808080809
$ws.Range("T2").Value = 608.8929292929294
$ws.Range("V2").Value = 392.2262626262627
$ws.Range("W2").Value = 313.4383838383839
$ws.Range("X2").Value = 260.9131313131313
$ws.Range("Y2").Value = 221.5191919191919
$ws.Range("B3").Value = 182.1252525252525
$ws.Range("C3").Value = 162.4282828282828
$ws.Range("D3").Value = 149.2969696969697
$ws.Range("E3").Value = 149.2969696969697
$ws.Range("F3").Value = 149.2969696969697
$ws.Range("I3").Value = 170.784
$ws.Range("J3").Value = 170.784
$ws.Range("K3").Value = 170.784
$ws.Range("L3").Value = 263.448
$ws.Range("M3").Value = 343.242
$ws.Range("N3").Value = 364.4805858585859
$ws.Range("O3").Value = 436.5525858585859
$ws.Range("P3").Value = 459.7185858585859
$ws.Range("Q3").Value = 485.4585858585859
$ws.Range("R3").Value = 485.4585858585859
$ws.Range("S3").Value = 464.4484848484849
$ws.Range("T3").Value = 333.1353535353535
$ws.Range("U3").Value = 333.1353535353535
$ws.Range("V3").Value = 333.1353535353535
$ws.Range("W3").Value = 254.3474747474747
$ws.Range("X3").Value = 254.3474747474747
$ws.Range("Y3").Value = 214.9535353535353
$ws.Range("B4").Value = 168.9939393939394
$ws.Range("C4").Value = 149.2969696969697
$ws.Range("D4").Value = 149.2969696969697
$ws.Range("E4").Value = 149.2969696969697
$ws.Range("F4").Value = 149.2969696969697
$ws.Range("J4").Value = 129.6
$ws.Range("K4").Value = 129.6
$ws.Range("L4").Value = 129.6
$ws.Range("M4").Value = 188.802
$ws.Range("N4").Value = 271.17
$ws.Range("O4").Value = 343.242
$ws.Range("P4").Value = 372.5292929292929
$ws.Range("T4").Value = 241.2161616161616
$ws.Range("U4").Value = 241.2161616161616
$ws.Range("V4").Value = 241.2161616161616
$ws.Range("W4").Value = 241.2161616161616
$ws.Range("X4").Value = 241.2161616161616

$ws = $wb.Worksheets.Item("Feed in from Type 5")
$ws.Range("S2").Value = 2.884000000000047
